$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 8.937933333333332
    "H2" = 26.8138
    "I2" = 0.2302024600837126
    "J2" = 0.2302024600837126
    "M2" = 0.9705896666666667
    "N2" = 2.911769
    "O2" = 0.02073452941466921
    "P2" = 0.02073452941466921
    "Q2" = 8.675065734688888
    "R2" = 78.07559161219999
    "S2" = 0.004773139679934955
    "T2" = 0.004773139679934954
    "G3" = 8.937933333333332
    "H3" = 26.8138
    "I3" = 0.2302024600837126
    "J3" = 0.2302024600837126
    "O3" = 0.5628689972673966
    "P3" = 0.5628689972673966
    "Q3" = 235.4972931219999
    "R3" = 2119.475638098
    "S3" = 0.1295738278758072
    "T3" = 0.1295738278758072
    "G4" = 8.937933333333332
    "H4" = 26.8138
    "I4" = 0.2302024600837126
    "J4" = 0.2302024600837126
    "M4" = 19.49164633333333
    "N4" = 58.47493899999999
    "O4" = 0.4163964733179342
    "P4" = 0.4163964733179341
    "Q4" = 174.2150354842444
    "R4" = 1567.9353193582
    "S4" = 0.09585549252797045
    "T4" = 0.09585549252797042
    "I5" = 0.5278886986241245
    "J5" = 0.5278886986241244
    "M5" = 0.9705896666666667
    "N5" = 2.911769
    "O5" = 0.02073452941466921
    "P5" = 0.02073452941466921
    "Q5" = 19.89322424920367
    "R5" = 179.039018242833
    "S5" = 0.01094552374929336
    "T5" = 0.01094552374929336
    "I6" = 0.5278886986241245
    "J6" = 0.5278886986241244
    "O6" = 0.5628689972673966
    "P6" = 0.5628689972673966
    "S6" = 0.2971321824633519
    "T6" = 0.2971321824633518
    "I7" = 0.5278886986241245
    "J7" = 0.5278886986241244
    "M7" = 19.49164633333333
    "N7" = 58.47493899999999
    "O7" = 0.4163964733179342
    "P7" = 0.4163964733179341
    "Q7" = 399.5011535892803
    "R7" = 3595.510382303522
    "S7" = 0.2198109924114792
    "T7" = 0.2198109924114792
    "G8" = 9.392449999999998
    "H8" = 28.17735
    "I8" = 0.241908841292163
    "J8" = 0.2419088412921629
    "M8" = 0.9705896666666667
    "N8" = 2.911769
    "O8" = 0.02073452941466921
    "P8" = 0.02073452941466921
    "Q8" = 9.116214914683331
    "R8" = 82.04593423214999
    "S8" = 0.005015865985440899
    "T8" = 0.005015865985440899
    "G9" = 9.392449999999998
    "H9" = 28.17735
    "I9" = 0.241908841292163
    "J9" = 0.2419088412921629
    "O9" = 0.5628689972673966
    "P9" = 0.5628689972673966
    "Q9" = 247.4729300714999
    "R9" = 2227.256370643499
    "S9" = 0.1361629869282376
    "T9" = 0.1361629869282375
    "G10" = 9.392449999999998
    "H10" = 28.17735
    "I10" = 0.241908841292163
    "J10" = 0.2419088412921629
    "M10" = 19.49164633333333
    "N10" = 58.47493899999999
    "O10" = 0.4163964733179342
    "P10" = 0.4163964733179341
    "Q10" = 183.0743136035166
    "R10" = 1647.66882243165
    "S10" = 0.1007299883784845
    "T10" = 0.1007299883784845
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
